$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. pitstop sheet: add tyre_before / tyre_after columns (H, I)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("pitstop")

$ws.Range("H1").Value = "tyre_before"
$ws.Range("I1").Value = "tyre_after"

$tyreData = @(
    @(3, 2),
    @(3, 2),
    @(3, 2),
    @(2, 3),
    @(3, 2),
    @(2, 2),
    @(2, 3),
    @(3, 2),
    @(3, 2),
    @(3, 2),
    @(3, 2),
    @(3, 2),
    @(2, 3),
    @(3, 2),
    @(2, 3),
    @(2, 3),
    @(3, 2),
    @(2, 3),
    @(3, 2),
    @(3, 3),
    @(2, 3),
    @(3, 2),
    @(2, 3),
    @(2, 3)
)

$row = 2
foreach ($pair in $tyreData) {
    $ws.Cells.Item($row, 8).Value = $pair[0]
    $ws.Cells.Item($row, 9).Value = $pair[1]
    $row = $row + 1
}

$ws.Activate()
$ws.Range("H25").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. New "weather" sheet (appended after the last existing sheet)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWeather = $wb.Worksheets.Add($null, $lastSheet)
$wsWeather.Name = "weather"

$wsWeather.Range("A1").Value = "Skycondition"
$wsWeather.Range("B1").Value = "Partly Cloudy"

$wsWeather.Range("A2").Value = "Precipation type"
$wsWeather.Range("B2").Value = "rain"

$wsWeather.Range("A3").Value = "Temperature"
$wsWeather.Range("B3").Value = "61.4°F"

$wsWeather.Range("A4").Value = "Humidity"
$wsWeather.Range("B4").Value = 0.62
$wsWeather.Range("B4").NumberFormat = "0%"

$wsWeather.Range("A5").Value = "Wind speed"
$wsWeather.Range("B5").Value = "7.86 mph"

$wsWeather.Range("A6").Value = "Wind bearing"
$wsWeather.Range("B6").Value = "329°"

$wsWeather.Columns.Item(1).AutoFit() | Out-Null
$wsWeather.Columns.Item(2).AutoFit() | Out-Null

$wsWeather.Activate()
$wsWeather.Range("D5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. New "altitude" sheet (appended after "weather", becomes active tab)
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAltitude = $wb.Worksheets.Add($null, $lastSheet2)
$wsAltitude.Name = "altitude"

$wsAltitude.Range("A1").Value = "delta"
$wsAltitude.Range("B1").Value = 102.2

$wsAltitude.Activate()
$wsAltitude.Range("B2").Select() | Out-Null
